$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.882.50"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.06%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.294.62"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.09%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.50"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.05"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.30%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.503"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.45%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.00"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +7.80%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0790"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.70%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "18.43"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +9.53%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.24%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.94"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.65%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.651.28"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.449.26"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +6.60%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.799"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.81%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.789.19"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.55"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +8.77%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0902"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.11"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.62%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.72%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.21"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.22"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +11.21%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.33%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.73%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.96"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.68%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +14.68%  "

$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.23"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.14%  "

$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.37"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.67%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.11"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.29%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.00"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.55%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +5.44%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.64"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.15%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.94%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0687"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.65%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.26%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.101"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.15%  "

$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.78"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.36%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.67%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.33%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0290"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.79%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.972.57"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.66%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.34%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.56"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.46%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.88"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.46%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.07"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.59%  "

$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.519.02"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.17%  "

$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.54"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.65"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.45%  "
